$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (municipio-nombre dimension) metadata rows
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"

# Column F (aragon dimension) metadata rows
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("F4").Value = "URI-Comunidad"

# Column G (sexo dimension) metadata rows
$ws.Range("G2").Value = "iaest-measure:sexo"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"

# Row 5 (mapping-aragon.xlsx / mapping-sexo.xlsx) is no longer needed
$ws.Range("F5:G5").Delete()
